$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update threshold values in column F for rows 8,9,10,13,14,15 (29 -> 19)
$ws.Range("F8").Value = 19
$ws.Range("F9").Value = 19
$ws.Range("F10").Value = 19
$ws.Range("F13").Value = 19
$ws.Range("F14").Value = 19
$ws.Range("F15").Value = 19

# Update the active selection to F22 as recorded in the saved view state
$ws.Range("F22").Select()
